$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update odds values in row 4 per the diff
$ws.Range("G4").Value = 2.7
$ws.Range("I4").Value = 2.7
$ws.Range("J4").Value = 3.25
$ws.Range("M4").Value = 1.07
$ws.Range("N4").Value = 9
$ws.Range("O4").Value = 1.33
$ws.Range("P4").Value = 3.25
$ws.Range("Q4").Value = 2.08
$ws.Range("R4").Value = 1.73
$ws.Range("S4").Value = 3.5
$ws.Range("T4").Value = 1.29
$ws.Range("U4").Value = 1.44
$ws.Range("V4").Value = 2.63
$ws.Range("AA4").Value = 11
$ws.Range("AD4").Value = 34
$ws.Range("AE4").Value = 9
$ws.Range("AJ4").Value = 9
$ws.Range("AM4").Value = 26
$ws.Range("AO4").Value = 34
